# Insert a new "Poll" slide (a duplicate of the existing "Poll Q1 ..."
# slide, which already carries the right layout / QR-code picture /
# placeholder formatting) right before the "Saving and restoring
# registers" divider slide, and retitle the duplicate to just "Poll".
#
# Before: slide 57 = "Saving and restoring registers" (subtitle divider)
# After:  slide 57 = new "Poll" slide
#         slide 58 = "Saving and restoring registers" (shifted down by one)
#         ... everything after shifts down by one slide as well.

$p = $ppt.ActivePresentation

# Slide 34 ("Poll Q1: Where do we allocate global variables") already has
# the title/body/picture layout we want to reuse for the new poll slide.
$template = $p.Slides.Item(34)

# Duplicate() inserts the copy immediately after the source slide and
# returns the SlideRange containing it.
$dupRange = $template.Duplicate()
$newSlide = $dupRange.Item(1)

# Move the duplicate into its final position: right before the old
# slide 57 ("Saving and restoring registers").
$newSlide.MoveTo(57)

# Retitle the duplicate from "Poll Q1: ..." to plain "Poll".
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Poll"
